$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list in column D
$ws.Range("D29").Value = 213.15
$ws.Range("D30").Value = 213.15
$ws.Range("D31").Value = 213.15
$ws.Range("D32").Value = 213.15
$ws.Range("D33").Value = 229.1
$ws.Range("D34").Value = 229.1
$ws.Range("D35").Value = 229.1
$ws.Range("D36").Value = 229.1
$ws.Range("D37").Value = 263
$ws.Range("D38").Value = 263
